$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 118.35
$ws.Range("I11").Value = 118.35
$ws.Range("K11").Value = 118.35
$ws.Range("M11").Value = 21.65000000000001
$ws.Range("H137").Value = 9903.519
$ws.Range("I137").Value = 1763.091
$ws.Range("K137").Value = 5289.272999999999
$ws.Range("M137").Value = -2739.272999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 62
$ws.Range("N5").ClearContents()
$ws.Range("H74").Value = 29407.867
$ws.Range("I74").Value = 2384.5715
$ws.Range("K74").Value = 2384.5715
$ws.Range("M74").Value = -1510.5715
$ws.Range("H77").Value = 29407.867
$ws.Range("I77").Value = 2384.5715
$ws.Range("K77").Value = 11922.8575
$ws.Range("M77").Value = -7554.8575
$ws.Range("H97").Value = 701.5454999999999
$ws.Range("I97").Value = 762.3333
$ws.Range("J97").Value = 628.6
$ws.Range("K97").Value = 762.3333
$ws.Range("L97").Value = 628.6
$ws.Range("M97").Value = -266.3333
$ws.Range("N97").Value = -1620.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 65
$ws.Range("N4").ClearContents()
$ws.Range("H15").Value = 480.66666
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H20").Value = 33849.395
$ws.Range("I20").Value = 18833.223
$ws.Range("J20").Value = 40962.316
$ws.Range("K20").Value = 18833.223
$ws.Range("L20").Value = 40962.316
$ws.Range("M20").Value = -18586.223
$ws.Range("N20").Value = -41456.316
$ws.Range("H99").Value = 19097.584
$ws.Range("I99").Value = 18168.422
$ws.Range("K99").Value = 18168.422
$ws.Range("M99").Value = -16670.422
$ws.Range("H134").Value = 150918.25
$ws.Range("J134").Value = 34446.668
$ws.Range("L134").Value = 103340.004
$ws.Range("N134").Value = -108410.004

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2118.6
$ws.Range("J2").Value = 2750
$ws.Range("L2").Value = 2750
$ws.Range("N2").Value = -2976
$ws.Range("H31").Value = 6831.026
$ws.Range("I31").Value = 1214.5264
$ws.Range("J31").Value = 12166.7
$ws.Range("K31").Value = 1214.5264
$ws.Range("L31").Value = 12166.7
$ws.Range("M31").Value = -919.5264
$ws.Range("N31").Value = -12756.7
$ws.Range("H34").Value = 6831.026
$ws.Range("I34").Value = 1214.5264
$ws.Range("J34").Value = 12166.7
$ws.Range("K34").Value = 1214.5264
$ws.Range("L34").Value = 12166.7
$ws.Range("M34").Value = -1012.5264
$ws.Range("N34").Value = -12570.7
$ws.Range("H59").Value = 28243.8
$ws.Range("I59").Value = 9104
$ws.Range("J59").Value = 33028.75
$ws.Range("K59").Value = 9104
$ws.Range("L59").Value = 33028.75
$ws.Range("M59").Value = -7959
$ws.Range("N59").Value = -35318.75
$ws.Range("H134").Value = 41673436
$ws.Range("I134").Value = 2111.9412
$ws.Range("J134").Value = 142875230
$ws.Range("K134").Value = 6335.823600000001
$ws.Range("L134").Value = 428625690
$ws.Range("M134").Value = -3800.823600000001
$ws.Range("N134").Value = -428630760

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 103.72222
$ws.Range("I2").Value = 143.33333
$ws.Range("J2").Value = 24.5
$ws.Range("K2").Value = 859.9999799999999
$ws.Range("L2").Value = 147
$ws.Range("M2").Value = -746.9999799999999
$ws.Range("N2").Value = -373
$ws.Range("H22").Value = 2782.3
$ws.Range("J22").Value = 3052.5557
$ws.Range("L22").Value = 9157.667099999999
$ws.Range("N22").Value = -9495.667099999999
$ws.Range("H27").Value = 2782.3
$ws.Range("J27").Value = 3052.5557
$ws.Range("L27").Value = 9157.667099999999
$ws.Range("N27").Value = -9361.667099999999
$ws.Range("H31").Value = 4577.5557
$ws.Range("J31").Value = 4719.8
$ws.Range("L31").Value = 14159.4
$ws.Range("N31").Value = -14735.4
$ws.Range("H63").Value = 2670.6667
$ws.Range("I63").Value = 2006
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 6018
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -5269
$ws.Range("N63").Value = -13498
$ws.Range("H66").Value = 2670.6667
$ws.Range("I66").Value = 2006
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 18054
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -14310
$ws.Range("N66").Value = -43488
$ws.Range("H68").Value = 997.087
$ws.Range("I68").Value = 396.5
$ws.Range("J68").Value = 1209.0588
$ws.Range("K68").Value = 1189.5
$ws.Range("L68").Value = 3627.1764
$ws.Range("M68").Value = -378.5
$ws.Range("N68").Value = -5249.1764
$ws.Range("H71").Value = 997.087
$ws.Range("I71").Value = 396.5
$ws.Range("J71").Value = 1209.0588
$ws.Range("K71").Value = 3568.5
$ws.Range("L71").Value = 10881.5292
$ws.Range("M71").Value = 487.5
$ws.Range("N71").Value = -18993.5292
$ws.Range("H98").Value = 895
$ws.Range("J98").Value = 895
$ws.Range("L98").Value = 2685
$ws.Range("N98").Value = -5681
$ws.Range("H124").Value = 5155.222
$ws.Range("I124").Value = 5049.625
$ws.Range("K124").Value = 15148.875
$ws.Range("M124").Value = -10238.875
$ws.Range("H131").Value = 1443.97
$ws.Range("I131").Value = 833
$ws.Range("J131").Value = 1462.866
$ws.Range("K131").Value = 2499
$ws.Range("L131").Value = 4388.598
$ws.Range("M131").Value = 2541
$ws.Range("N131").Value = -14468.598
$ws.Range("H134").Value = 6194.838
$ws.Range("I134").Value = 2670.9
$ws.Range("K134").Value = 8012.700000000001
$ws.Range("M134").Value = -2942.700000000001
$ws.Range("H137").Value = 4999
$ws.Range("I137").Value = 3760.1428
$ws.Range("J137").Value = 6083
$ws.Range("K137").Value = 11280.4284
$ws.Range("L137").Value = 18249
$ws.Range("M137").Value = -6180.428400000001
$ws.Range("N137").Value = -28449

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 20000
$ws.Range("J53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("N53").Value = -21262
$ws.Range("H70").Value = 17499.5
$ws.Range("J70").Value = 17499.5
$ws.Range("L70").Value = 17499.5
$ws.Range("N70").Value = -18039.5
$ws.Range("H73").Value = 17499.5
$ws.Range("J73").Value = 17499.5
$ws.Range("L73").Value = 17499.5
$ws.Range("N73").Value = -19371.5
$ws.Range("H97").Value = 1508.6471
$ws.Range("I97").Value = 1419.2307
$ws.Range("K97").Value = 1419.2307
$ws.Range("M97").Value = -923.2307000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1050.6875
$ws.Range("J16").Value = 629.6667
$ws.Range("L16").Value = 629.6667
$ws.Range("N16").Value = -969.6667
$ws.Range("H22").Value = 742.25
$ws.Range("I22").Value = 689.6667
$ws.Range("K22").Value = 689.6667
$ws.Range("M22").Value = -394.6667
$ws.Range("H27").Value = 742.25
$ws.Range("I27").Value = 689.6667
$ws.Range("K27").Value = 689.6667
$ws.Range("M27").Value = -582.6667
$ws.Range("H93").Value = 9299.286
$ws.Range("I93").Value = 14050
$ws.Range("J93").Value = 2965
$ws.Range("K93").Value = 14050
$ws.Range("L93").Value = 2965
$ws.Range("M93").Value = -12802
$ws.Range("N93").Value = -5461

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10800.637
$ws.Range("I62").Value = 9513.799999999999
$ws.Range("J62").Value = 11873
$ws.Range("K62").Value = 9513.799999999999
$ws.Range("L62").Value = 11873
$ws.Range("M62").Value = -8889.799999999999
$ws.Range("N62").Value = -13121
$ws.Range("H65").Value = 10800.637
$ws.Range("I65").Value = 9513.799999999999
$ws.Range("J65").Value = 11873
$ws.Range("K65").Value = 47569
$ws.Range("L65").Value = 59365
$ws.Range("M65").Value = -44449
$ws.Range("N65").Value = -65605
$ws.Range("H107").Value = 830.69696
$ws.Range("I107").Value = 894.7037
$ws.Range("K107").Value = 2684.1111
$ws.Range("M107").Value = -764.1111000000001
